$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date in A1 (step forward one day, from 45308 to 45309)
$ws.Range("A1").Value = 45309

# Update prices in "Para CARGA" section (rows 23-27)
$ws.Range("D23").Value = 1931
$ws.Range("D24").Value = 2106
$ws.Range("D25").Value = 2250
$ws.Range("D26").Value = 2320
$ws.Range("D27").Value = 2790

# Update prices in "Para DESCARGA" section (rows 34-35)
$ws.Range("D34").Value = 1351
$ws.Range("D35").Value = 1651
